$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style from an existing header cell (A1) to AD1:AF1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Data rows 2-52
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 101   # AD
    $ws.Cells.Item($r, 31).Value = 61    # AE
    $ws.Cells.Item($r, 32).Value = 0     # AF
}
